$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN['German Shorthaired Pointer'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
         coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newQuery

# The query got two shorter (two RETURN lines removed), so the wrapped row
# shrinks from its old auto-computed height.
$ws.Rows.Item(4).RowHeight = 217.5

# Update the view: active cell moves from C5 to B4, and the top-left cell reset (scrolled to A1)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
